$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# New explanatory text for the "Описание" sheet (sheet1).
# Values are written in the order that reproduces the shared-string table
# indices used by the target workbook (35..42).
$s35 = "В файле две вкладки. Первая — справочная. В соответствии с указанием в ней будут выставлены настройки массовой загрузки."
$s36 = 'Вторая — слова и переводы."1" и "2" — это номера языков с первой вкладки.'
$s37 = "http://slovari.yandex.ru/%D0%BA%D0%BE%D1%88%D0%BA%D0%B0/%D0%BF%D0%B5%D1%80%D0%B5%D0%B2%D0%BE%D0%B4/"
$s38 = 'тут ограничением будет "(для лазания на столбы)" или "(электрического крана)".'
$s39 = '"Написание" — для слова заполняется, если есть официальная орфография. Для ижорского ВСЕГДА пусто.'
$s40 = '"Произношение" может быть, например, в "учебной", в МФА или просто в той, в которой приведено в источнике. Варианты произношения (или даже написания, такое может быть, если язык имеет несколько литературных норм, как английский) указываются через вертикальную черту. Диалект в скобках (везде одинаково должно быть, например ala и soi).'
$s41 = '"Часть речи" желательно указывать принятым английским сокращением, но можно и по-русски или по-фински. Главное, чтобы во всём файле было одинаково.'
$s42 = '"Ограничение перевода" — это, например:'

$ws1.Range("A9").Value = $s35
$ws1.Range("A10").Value = $s36
$ws1.Range("A15").Value = $s37
$ws1.Range("A16").Value = $s38
$ws1.Range("A11").Value = $s39
$ws1.Range("A12").Value = $s40
$ws1.Range("A13").Value = $s41
$ws1.Range("A14").Value = $s42

# Turn the URL in A15 into a real hyperlink.
$ws1.Hyperlinks.Add($ws1.Range("A15"), $s37) | Out-Null

# Update the selections: "Слова" (sheet2) keeps a plain selection but is no
# longer the active tab, while "Описание" (sheet1) becomes the active tab
# with C19 selected.
$ws2.Range("A42").Select() | Out-Null
$ws1.Range("C19").Select() | Out-Null
